$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 302, shifting existing rows 302-338 down to 303-339.
$ws.Rows.Item(302).Insert()

# Populate the new row 302 with this week's entry (same template as the
# surrounding Jengibre / Vega Modelo de Temuco rows, with the new date,
# volume and price figures).
$ws.Cells.Item(302, 1).Value = 10
$ws.Cells.Item(302, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(302, 3).Value = "La Araucanía"
$ws.Cells.Item(302, 4).Value = 45142
$ws.Cells.Item(302, 5).Value = 9
$ws.Cells.Item(302, 6).Value = 100114007
$ws.Cells.Item(302, 7).Value = "Jengibre"
$ws.Cells.Item(302, 8).Value = "Sin especificar"
$ws.Cells.Item(302, 9).Value = "Primera"
$ws.Cells.Item(302, 10).Value = 120
$ws.Cells.Item(302, 11).Value = 24000
$ws.Cells.Item(302, 12).Value = 24000
$ws.Cells.Item(302, 13).Value = 24000
$ws.Cells.Item(302, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(302, 15).Value = "Perú"
$ws.Cells.Item(302, 16).Value = 1846
$ws.Cells.Item(302, 17).Value = 13
$ws.Cells.Item(302, 18).Value = "Hortaliza"
